$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 3819.8
$ws.Range("I13").Value = 2250
$ws.Range("J13").Value = 4866.3335
$ws.Range("K13").Value = 2250
$ws.Range("L13").Value = 4866.3335
$ws.Range("M13").Value = -2081
$ws.Range("N13").Value = -5204.3335

$ws.Range("H138").Value = 2973.7812
$ws.Range("J138").Value = 3343.6938
$ws.Range("L138").Value = 10031.0814
$ws.Range("N138").Value = -20311.0814

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2737
$ws.Range("I2").Value = 2945.7778
$ws.Range("K2").Value = 2945.7778
$ws.Range("M2").Value = -2832.7778

$ws.Range("H32").Value = 5765.289
$ws.Range("I32").Value = 3739.95
$ws.Range("K32").Value = 3739.95
$ws.Range("M32").Value = -3452.95

$ws.Range("H61").Value = 49913.6
$ws.Range("I61").Value = 2744.2856
$ws.Range("J61").Value = 75312.46000000001
$ws.Range("K61").Value = 2744.2856
$ws.Range("L61").Value = 75312.46000000001
$ws.Range("M61").Value = -2532.2856
$ws.Range("N61").Value = -75736.46000000001

$ws.Range("H74").Value = 362302.4
$ws.Range("I74").Value = 900000
$ws.Range("J74").Value = 227878
$ws.Range("K74").Value = 900000
$ws.Range("L74").Value = 227878
$ws.Range("M74").Value = -899126
$ws.Range("N74").Value = -229626

$ws.Range("H77").Value = 362302.4
$ws.Range("I77").Value = 900000
$ws.Range("J77").Value = 227878
$ws.Range("K77").Value = 4500000
$ws.Range("L77").Value = 1139390
$ws.Range("M77").Value = -4495632
$ws.Range("N77").Value = -1148126

$ws.Range("H116").Value = 2737
$ws.Range("I116").Value = 2945.7778
$ws.Range("K116").Value = 2945.7778
$ws.Range("M116").Value = -651.7777999999998

$ws.Range("H122").Value = 304833.34
$ws.Range("J122").Value = 304833.34
$ws.Range("L122").Value = 914500.02
$ws.Range("N122").Value = -919400.02

$ws.Range("H132").Value = 5925.4355
$ws.Range("I132").Value = 6176.607
$ws.Range("J132").Value = 3581.1667
$ws.Range("K132").Value = 18529.821
$ws.Range("L132").Value = 10743.5001
$ws.Range("M132").Value = -15999.821
$ws.Range("N132").Value = -15803.5001

$ws.Range("H136").Value = 49913.6
$ws.Range("I136").Value = 2744.2856
$ws.Range("J136").Value = 75312.46000000001
$ws.Range("K136").Value = 8232.856800000001
$ws.Range("L136").Value = 225937.38
$ws.Range("M136").Value = -5682.856800000001
$ws.Range("N136").Value = -231037.38

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2737
$ws.Range("I3").Value = 2945.7778
$ws.Range("K3").Value = 2945.7778
$ws.Range("M3").Value = -2831.7778

$ws.Range("H76").Value = 43184.4
$ws.Range("J76").Value = 43184.4
$ws.Range("L76").Value = 43184.4
$ws.Range("N76").Value = -43814.4

$ws.Range("H79").Value = 43184.4
$ws.Range("J79").Value = 43184.4
$ws.Range("L79").Value = 43184.4
$ws.Range("N79").Value = -45368.4

$ws.Range("H94").Value = 48560.156
$ws.Range("I94").Value = 837.3570999999999
$ws.Range("K94").Value = 837.3570999999999
$ws.Range("M94").Value = -386.3570999999999

$ws.Range("H99").Value = 2166.0417
$ws.Range("I99").Value = 2008.1875
$ws.Range("K99").Value = 2008.1875
$ws.Range("M99").Value = -510.1875

$ws.Range("H134").Value = 3373.383
$ws.Range("I134").Value = 2282.359
$ws.Range("J134").Value = 8692.125
$ws.Range("K134").Value = 6847.076999999999
$ws.Range("L134").Value = 26076.375
$ws.Range("M134").Value = -4312.076999999999
$ws.Range("N134").Value = -31146.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1281.25
$ws.Range("I22").Value = 1129
$ws.Range("J22").Value = 1332
$ws.Range("K22").Value = 1129
$ws.Range("L22").Value = 1332
$ws.Range("M22").Value = -779
$ws.Range("N22").Value = -2032

$ws.Range("H58").Value = 3969.3948
$ws.Range("I58").Value = 3099.92
$ws.Range("J58").Value = 5641.4614
$ws.Range("K58").Value = 3099.92
$ws.Range("L58").Value = 5641.4614
$ws.Range("M58").Value = -2896.92
$ws.Range("N58").Value = -6047.4614

$ws.Range("H99").Value = 8278.857
$ws.Range("I99").Value = 5592.9
$ws.Range("K99").Value = 5592.9
$ws.Range("M99").Value = -4094.9

$ws.Range("H105").Value = 2020.1765
$ws.Range("I105").Value = 1371.0834
$ws.Range("J105").Value = 3578
$ws.Range("K105").Value = 1371.0834
$ws.Range("L105").Value = 3578
$ws.Range("M105").Value = 375.9166
$ws.Range("N105").Value = -7072

$ws.Range("H107").Value = 42198.043
$ws.Range("I107").Value = 71302.766
$ws.Range("K107").Value = 71302.766
$ws.Range("M107").Value = -69382.766

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = 0

$ws.Range("H126").Value = 8278.857
$ws.Range("I126").Value = 5592.9
$ws.Range("K126").Value = 16778.7
$ws.Range("M126").Value = -14308.7

$ws.Range("H132").Value = 3647.761
$ws.Range("I132").Value = 1691.475
$ws.Range("J132").Value = 16689.666
$ws.Range("K132").Value = 5074.424999999999
$ws.Range("L132").Value = 50068.99800000001
$ws.Range("M132").Value = -2544.424999999999
$ws.Range("N132").Value = -55128.99800000001

$ws.Range("H136").Value = 3969.3948
$ws.Range("I136").Value = 3099.92
$ws.Range("J136").Value = 5641.4614
$ws.Range("K136").Value = 9299.76
$ws.Range("L136").Value = 16924.3842
$ws.Range("M136").Value = -6749.76
$ws.Range("N136").Value = -22024.3842

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 82940.27
$ws.Range("I122").Value = 274.5
$ws.Range("K122").Value = 2470.5
$ws.Range("M122").Value = -20.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 319876
$ws.Range("I5").Value = 319876
$ws.Range("K5").Value = 319876
$ws.Range("M5").Value = -319764

$ws.Range("H13").Value = 5702.25
$ws.Range("J13").Value = 1399.5
$ws.Range("L13").Value = 1399.5
$ws.Range("N13").Value = -1677.5

$ws.Range("H102").Value = 37761.88
$ws.Range("I102").Value = 47654.695
$ws.Range("K102").Value = 47654.695
$ws.Range("M102").Value = -46032.695

$ws.Range("H122").Value = 154679.17
$ws.Range("I122").Value = 226268.75
$ws.Range("K122").Value = 678806.25
$ws.Range("M122").Value = -676356.25

$ws.Range("H126").Value = 86592.27
$ws.Range("J126").Value = 8671
$ws.Range("L126").Value = 26013
$ws.Range("N126").Value = -30953

$ws.Range("H132").Value = 6832.1665
$ws.Range("I132").Value = 6832.1665
$ws.Range("K132").Value = 20496.4995
$ws.Range("M132").Value = -17966.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1168.2354
$ws.Range("I61").Value = 944.6923
$ws.Range("J61").Value = 1894.75
$ws.Range("K61").Value = 944.6923
$ws.Range("L61").Value = 1894.75
$ws.Range("M61").Value = -742.6923
$ws.Range("N61").Value = -2298.75

$ws.Range("H93").Value = 1226.6
$ws.Range("I93").Value = 1091.9166
$ws.Range("J93").Value = 1765.3334
$ws.Range("K93").Value = 1091.9166
$ws.Range("L93").Value = 1765.3334
$ws.Range("M93").Value = 156.0834
$ws.Range("N93").Value = -4261.3334

$ws.Range("H113").Value = 1168.2354
$ws.Range("I113").Value = 944.6923
$ws.Range("J113").Value = 1894.75
$ws.Range("K113").Value = 944.6923
$ws.Range("L113").Value = 1894.75
$ws.Range("M113").Value = 1225.3077
$ws.Range("N113").Value = -6234.75

$ws.Range("H122").Value = 6296.3335
$ws.Range("I122").Value = 3458.375
$ws.Range("K122").Value = 10375.125
$ws.Range("M122").Value = -7925.125

$ws.Range("H136").Value = 5757.5454
$ws.Range("I136").Value = 3505
$ws.Range("K136").Value = 10515
$ws.Range("M136").Value = -7965

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 16961634
$ws.Range("I5").Value = 20314000
$ws.Range("J5").Value = 199800
$ws.Range("K5").Value = 20314000
$ws.Range("L5").Value = 199800
$ws.Range("M5").Value = -20313888
$ws.Range("N5").Value = -200024

$ws.Range("H11").Value = 15003000
$ws.Range("I11").Value = 20000000
$ws.Range("J11").Value = 13337333
$ws.Range("K11").Value = 20000000
$ws.Range("L11").Value = 13337333
$ws.Range("M11").Value = -19999858
$ws.Range("N11").Value = -13337617

$ws.Range("H13").Value = 3566.6667
$ws.Range("I13").Value = 2850
$ws.Range("J13").Value = 5000
$ws.Range("K13").Value = 2850
$ws.Range("L13").Value = 5000
$ws.Range("M13").Value = -2710
$ws.Range("N13").Value = -5280

$ws.Range("H122").Value = 5481.5186
$ws.Range("I122").Value = 5183.5
$ws.Range("J122").Value = 6792.8
$ws.Range("K122").Value = 15550.5
$ws.Range("L122").Value = 20378.4
$ws.Range("M122").Value = -13100.5
$ws.Range("N122").Value = -25278.4

$ws.Range("H132").Value = 3120.4424
$ws.Range("I132").Value = 1824.9318
$ws.Range("K132").Value = 5474.7954
$ws.Range("M132").Value = -2944.7954

$ws.Range("H136").Value = 9969.883
$ws.Range("I136").Value = 11441.5
$ws.Range("J136").Value = 6438
$ws.Range("K136").Value = 34324.5
$ws.Range("L136").Value = 19314
$ws.Range("M136").Value = -31774.5
$ws.Range("N136").Value = -24414
